$wb = $excel.ActiveWorkbook

# Worksheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7595
$ws.Range("F4").Value = 73
$ws.Range("F5").Value = 4778
$ws.Range("F10").Value = 134
$ws.Range("F12").Value = 758
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 66
$ws.Range("F15").Value = 252
$ws.Range("F16").Value = 15
$ws.Range("F17").Value = 248
$ws.Range("F19").Value = 381
$ws.Range("F20").Value = 138
$ws.Range("F21").Value = 1080
$ws.Range("F23").Value = 566
$ws.Range("F24").Value = 2156
$ws.Range("F25").Value = 688
$ws.Range("F26").Value = 39
$ws.Range("F27").Value = 37
$ws.Range("F29").Value = 593
$ws.Range("F30").Value = 39

# Worksheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 429

# Worksheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 429
$ws.Range("F3").Value = 7595
$ws.Range("F5").Value = 73
$ws.Range("F7").Value = 4780
$ws.Range("F13").Value = 134
$ws.Range("F18").Value = 758
$ws.Range("F19").Value = 27
$ws.Range("F20").Value = 66
$ws.Range("F21").Value = 252
$ws.Range("F23").Value = 15
$ws.Range("F26").Value = 248
$ws.Range("F28").Value = 381
$ws.Range("F29").Value = 138
$ws.Range("F30").Value = 1080
$ws.Range("F32").Value = 566
$ws.Range("F33").Value = 2156
$ws.Range("F34").Value = 688
$ws.Range("F35").Value = 39
$ws.Range("F36").Value = 37
$ws.Range("F38").Value = 593
$ws.Range("F39").Value = 39
